$d = $word.ActiveDocument

# The "Requisitos" bullet paragraph lists three course-requirement lines,
# each its own run terminated with a manual line break (<w:br/>). Move the
# "LOB1012 -  Estatística  (Requisito fraco)" line from the end of the
# list to the front, ahead of "LOQ4095 -  Química Geral Experimental
# (Requisito fraco)" / "LOQ4098 -  Fundamentos de Química para
# Engenharia II (Requisito fraco)".

# Locate the target paragraph robustly (rather than assuming it is the
# document's last paragraph): the one whose text contains the LOQ4095
# requirement line.
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*LOQ4095*") {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not locate the Requisitos paragraph"
}

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>' +
       '<w:r><w:t>LOB1012 -  Estatística  (Requisito fraco)</w:t><w:br/></w:r>' +
       '<w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito fraco)</w:t><w:br/></w:r>' +
       '<w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)</w:t><w:br/></w:r>' +
       '</w:p>'

$target.Range.InsertXML($xml) | Out-Null
